$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the new centered (horizontal+vertical) cell style in one shot ---
# Using a scratch cell + Copy/PasteSpecial(formats) avoids leaving an orphan
# "horizontal-only" style xf behind (each direct property write mints its own xf).
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("A2:P8").PasteSpecial(-4122)
$ws.Range("Z1").Clear()
# F4 has no cell at all in the target sheet (never repopulated after the
# original data loss) - drop the pasted-format placeholder entirely.
$ws.Range("F4").Clear()

# --- Populate the rebuilt data rows, column-by-column (matches the shared
# string insertion order baked into the target workbook) ---

# Column A - Flight Type
$ws.Range("A2").Value = "Oneway"
$ws.Range("A3").Value = "Return"
$ws.Range("A4").Value = "Multi"
$ws.Range("A5").Value = "Oneway"
$ws.Range("A6").Value = "Return"
$ws.Range("A7").Value = "Multi"
$ws.Range("A8").Value = ""

# Column B - Extra Paths
$ws.Range("B2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = ""

# Column C - From
$ws.Range("C2").Value = "Melbourne"
$ws.Range("C3").Value = "Melbourne"
$ws.Range("C4").Value = "Melbourne"
$ws.Range("C5").Value = "Melbourne"
$ws.Range("C6").Value = "Sydney"
$ws.Range("C7").Value = "Sydney"
$ws.Range("C8").Value = ""

# Column D - Destination
$ws.Range("D2").Value = "Sydney"
$ws.Range("D3").Value = "Tasmania"
$ws.Range("D4").Value = "Adelaide"
$ws.Range("D5").Value = "Adelaide"
$ws.Range("D6").Value = "Melbourne"
$ws.Range("D7").Value = "Hobart"
$ws.Range("D8").Value = ""

# Column E - Departure Date
$ws.Range("E2").Value = 25012022
$ws.Range("E3").Value = 25022022
$ws.Range("E4").Value = 25012022
$ws.Range("E5").Value = 25012022
$ws.Range("E6").Value = 25022022
$ws.Range("E7").Value = 25012022
$ws.Range("E8").Value = ""

# Column F - header relabelled "Arrival Date" -> "Return Date"
$ws.Range("F1").Value = "Return Date"
$ws.Range("F2").Value = ""
$ws.Range("F3").Value = 27022022
# F4 intentionally left untouched (no cell at all in the target sheet)
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = 27022022
$ws.Range("F7").Value = ""
$ws.Range("F8").Value = ""

# Column G - Adults
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = ""

# Column H - Children
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 4
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 2
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = ""

# Column I - Infants
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 2
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = ""

# Column J - Class
$ws.Range("J2").Value = "Economy"
$ws.Range("J3").Value = "First Class"
$ws.Range("J4").Value = "Business"
$ws.Range("J5").Value = "Premium Economy"
$ws.Range("J6").Value = "Economy"
$ws.Range("J7").Value = "Economy"
$ws.Range("J8").Value = ""

# Column K - Multi Destination
$ws.Range("K2").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("K4").Value = "Sydney"
$ws.Range("K5").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("K7").Value = "Melbourne"
$ws.Range("K8").Value = ""

# Column L - Multi Destination Date
$ws.Range("L2").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("L4").Value = 27012022
$ws.Range("L5").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("L7").Value = 27012022
$ws.Range("L8").Value = ""

# Column M - Multi Destination 2
$ws.Range("M2").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("M4").Value = "Melbourne"
$ws.Range("M5").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("M8").Value = ""

# Column N - Multi Destination Date 2
$ws.Range("N2").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("N4").Value = 28012022
$ws.Range("N5").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("N8").Value = ""

# Column O - Multi Destination 3
$ws.Range("O2").Value = ""
$ws.Range("O3").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("O7").Value = ""
$ws.Range("O8").Value = ""

# Column P - Multi Destination Date 3
$ws.Range("P2").Value = ""
$ws.Range("P3").Value = ""
$ws.Range("P4").Value = ""
$ws.Range("P5").Value = ""
$ws.Range("P6").Value = ""
$ws.Range("P7").Value = ""
$ws.Range("P8").Value = ""

# --- Column J is now wider (holds "Premium Economy" etc.) and no longer
# relies on autofit/bestFit, just a fixed width of 19 characters ---
$ws.Range("J1").ColumnWidth = 18.14

# --- Sheet-level view bits ---
$ws.Range("A2:P7").Select()
